$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILESTREAM")

# Insert a new row above row 2, shifting existing rows (and their formatting) down.
$ws.Rows.Item(2).Insert()

# Populate the new row with the BPQ deployment/DA entry.
$ws.Range("A2").Value = "NHANES-2015-2016-BPQ_I"
$ws.Range("B2").Value = "NHANES-BPQ"
$ws.Range("C2").Value = "nhanes-kb:DPL-BPQ_I-QUESTIONNAIRE"
$ws.Range("E2").Value = "example@example.com"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:example@example.com") | Out-Null
$ws.Range("F2").Value = "Public"

# Match the new selection/active cell recorded in the workbook view.
$ws.Range("A2").Select()

$wb.Save()
